# financial_worksheets_module-5.xlsx — apply commit "Add find wiki. Change png
# and pdf file permissions" (the message is boilerplate; the real edit is to
# the workbook content as captured by the OOXML diff):
#   1. Fill in the learner's name on the "Comprehensive Budget" sheet (the
#      other three sheets already show it) -> this removes the now-unused
#      "[Enter your name here]" shared string when the file is saved.
#   2. Change the Savings Goal's time horizon from 60 to 84 months, which
#      ripples into the recalculated monthly-savings-needed figure.
#   3. Fill in "Budgeted" (column E) figures on the Comprehensive Budget
#      sheet for rows 6-60 (totals in rows 12/19/40/55/61/62 recalc from
#      these automatically).
#   4. Leave the workbook with the same selections/active sheet the author
#      ended up with.

$wb = $excel.ActiveWorkbook

$wsBalance = $wb.Worksheets.Item("Balance Sheet")
$wsCashFlow = $wb.Worksheets.Item("Cash Flow")
$wsSavingsGoal = $wb.Worksheets.Item("Savings Goal")
$wsBudget = $wb.Worksheets.Item("Comprehensive Budget")

# 1) Name the budget on the Comprehensive Budget sheet, matching the other
#    tabs ("Mark Lucernas").
$wsBudget.Range("B2").Value = "Mark Lucernas"

# 2) Savings Goal: number of months to reach the goal 60 -> 84 (the
#    dependent "Monthly Savings Needed to Reach Goal" formula recalculates).
$wsSavingsGoal.Range("C9").Value = 84

# 3) Comprehensive Budget: "Budgeted" column (E) entries for rows 6-60.
$wsBudget.Range("E6").Value = 627
$wsBudget.Range("E7").Value = 0
$wsBudget.Range("E8").Value = 0
$wsBudget.Range("E9").Value = 0
$wsBudget.Range("E10").Value = 0
$wsBudget.Range("E11").Value = 0

$wsBudget.Range("E13").Value = 0
$wsBudget.Range("E14").Value = 0
$wsBudget.Range("E15").Value = 0
$wsBudget.Range("E16").Value = 0
$wsBudget.Range("E17").Value = 0
$wsBudget.Range("E18").Value = 0

$wsBudget.Range("E20").Value = 0
$wsBudget.Range("E21").Value = 33
$wsBudget.Range("E22").Value = 40
$wsBudget.Range("E23").Value = 0
$wsBudget.Range("E24").Value = 33
$wsBudget.Range("E25").Value = 0
$wsBudget.Range("E26").Value = 30
$wsBudget.Range("E27").Value = 0
$wsBudget.Range("E28").Value = 0
$wsBudget.Range("E29").Value = 0
$wsBudget.Range("E30").Value = 120
$wsBudget.Range("E31").Value = 0
$wsBudget.Range("E32").Value = 0
$wsBudget.Range("E33").Value = 0
$wsBudget.Range("E34").Value = 0
$wsBudget.Range("E35").Value = 0
$wsBudget.Range("E36").Value = 0
$wsBudget.Range("E37").Value = 0
$wsBudget.Range("E38").Value = 0
$wsBudget.Range("E39").Value = 0

$wsBudget.Range("E41").Value = 20
$wsBudget.Range("E42").Value = 0
$wsBudget.Range("E43").Value = 0
$wsBudget.Range("E44").Value = 0
$wsBudget.Range("E45").Value = 20
$wsBudget.Range("E46").Value = 0
$wsBudget.Range("E47").Value = 0
$wsBudget.Range("E48").Value = 0
$wsBudget.Range("E49").Value = 13
$wsBudget.Range("E50").Value = 0
$wsBudget.Range("E51").Value = 0
$wsBudget.Range("E52").Value = 0
$wsBudget.Range("E53").Value = 0
$wsBudget.Range("E54").Value = 0

$wsBudget.Range("E56").Value = 50
$wsBudget.Range("E57").Value = 0
$wsBudget.Range("E58").Value = 50
$wsBudget.Range("E59").Value = 20
$wsBudget.Range("E60").Value = 200

# 4) Restore per-sheet selections (each ends up matching the author's final
#    cursor position on that tab), finishing on "Comprehensive Budget" so it
#    is the workbook's active tab on save.
[void]$wsBalance.Range("E29").Select()

[void]$wsCashFlow.Activate()
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 7
[void]$wsCashFlow.Range("R31").Select()

[void]$wsSavingsGoal.Range("B2").Select()

[void]$wsBudget.Activate()
[void]$wsBudget.Range("I56").Select()
